# Atualização de bases das ligas, do dia: 03-03-2024 às 00:35
# Refresh match-odds/result data for "Uruguay Primera División" sheet.
#
# Rows 114/115 and 117-120 had their underlying match records re-ordered
# (the scraped source re-sorted matches with duplicate/near-duplicate
# timestamps), so every data column (id + F..AC) for those rows is
# rewritten to the record that now belongs there. Rows 138-142 are
# future fixtures that have since been played, so results (H/I/J) and
# odds/PL columns (N..AC) are updated/added with fresh values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 114
$ws.Range("B114").Value = 7559468
$ws.Range("F114").Value = "Liverpool Montevideo"
$ws.Range("G114").Value = "CA River Plate"
$ws.Range("H114").Value = 2
$ws.Range("I114").Value = 1
$ws.Range("J114").Value = "H"
$ws.Range("K114").Value = 1.7
$ws.Range("L114").Value = 3
$ws.Range("M114").Value = 5.75
$ws.Range("N114").Value = 1.833
$ws.Range("P114").Value = 4.5
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 1.925
$ws.Range("S114").Value = 1.925
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 2.025
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 0.833
$ws.Range("X114").Value = -1
$ws.Range("Z114").Value = 0.925
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = 1.025
$ws.Range("AC114").Value = -1
# Row 115
$ws.Range("B115").Value = 7559469
$ws.Range("F115").Value = "Montevideo Wanderers"
$ws.Range("G115").Value = "Penarol"
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 4.75
$ws.Range("L115").Value = 3.4
$ws.Range("M115").Value = 1.7
$ws.Range("N115").Value = 2.7
$ws.Range("P115").Value = 2.45
$ws.Range("Q115").Value = 0
$ws.Range("R115").Value = 2.05
$ws.Range("S115").Value = 1.8
$ws.Range("T115").Value = 2.5
$ws.Range("U115").Value = 1.975
$ws.Range("V115").Value = 1.875
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 2.2
$ws.Range("Z115").Value = 0
$ws.Range("AA115").Value = -0
$ws.Range("AB115").Value = -1
$ws.Range("AC115").Value = 0.875
# Row 117
$ws.Range("B117").Value = 7013885
$ws.Range("F117").Value = "La Luz"
$ws.Range("G117").Value = "Atletico Fenix Montevideo"
$ws.Range("K117").Value = 3
$ws.Range("L117").Value = 3
$ws.Range("M117").Value = 2.4
$ws.Range("N117").Value = 2.9
$ws.Range("O117").Value = 2.75
$ws.Range("P117").Value = 2.6
$ws.Range("Q117").Value = 0
$ws.Range("R117").Value = 2.025
$ws.Range("S117").Value = 1.825
$ws.Range("T117").Value = 2
$ws.Range("U117").Value = 2.025
$ws.Range("V117").Value = 1.825
$ws.Range("Y117").Value = 1.6
$ws.Range("AA117").Value = 0.825
$ws.Range("AB117").Value = 0
$ws.Range("AC117").Value = -0
# Row 118
$ws.Range("B118").Value = 7013702
$ws.Range("F118").Value = "Defensor Sporting"
$ws.Range("G118").Value = "Danubio"
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 2
$ws.Range("J118").Value = "A"
$ws.Range("K118").Value = 1.8
$ws.Range("L118").Value = 3.6
$ws.Range("M118").Value = 4.2
$ws.Range("N118").Value = 1.8
$ws.Range("O118").Value = 3.6
$ws.Range("P118").Value = 4.2
$ws.Range("R118").Value = 2.05
$ws.Range("S118").Value = 1.8
$ws.Range("T118").Value = 2.25
$ws.Range("U118").Value = 1.85
$ws.Range("V118").Value = 2
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = 3.2
$ws.Range("AA118").Value = 0.8
$ws.Range("AB118").Value = -0.5
$ws.Range("AC118").Value = 0.5
# Row 119
$ws.Range("B119").Value = 7013409
$ws.Range("F119").Value = "Nacional De Football"
$ws.Range("G119").Value = "Torque"
$ws.Range("H119").Value = 1
$ws.Range("J119").Value = "D"
$ws.Range("K119").Value = 1.666
$ws.Range("L119").Value = 3.9
$ws.Range("M119").Value = 4.5
$ws.Range("N119").Value = 1.615
$ws.Range("O119").Value = 4
$ws.Range("P119").Value = 4.75
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 1.8
$ws.Range("S119").Value = 2.05
$ws.Range("T119").Value = 2.75
$ws.Range("U119").Value = 1.95
$ws.Range("V119").Value = 1.9
$ws.Range("X119").Value = 3
$ws.Range("Y119").Value = -1
$ws.Range("AA119").Value = 1.05
$ws.Range("AC119").Value = 0.8999999999999999
# Row 120
$ws.Range("B120").Value = 7013886
$ws.Range("F120").Value = "Racing Club de Montevideo"
$ws.Range("G120").Value = "Cerro"
$ws.Range("I120").Value = 1
$ws.Range("K120").Value = 2.25
$ws.Range("L120").Value = 3.1
$ws.Range("M120").Value = 3.25
$ws.Range("N120").Value = 2.25
$ws.Range("O120").Value = 2.875
$ws.Range("P120").Value = 3.5
$ws.Range("Q120").Value = -0.25
$ws.Range("R120").Value = 1.95
$ws.Range("S120").Value = 1.9
$ws.Range("U120").Value = 1.925
$ws.Range("V120").Value = 1.925
$ws.Range("Y120").Value = 2.5
$ws.Range("AA120").Value = 0.8999999999999999
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = 0.925
# Row 138
$ws.Range("H138").Value = 1
$ws.Range("I138").Value = 2
$ws.Range("J138").Value = "A"
$ws.Range("N138").Value = 3.6
$ws.Range("P138").Value = 2.05
$ws.Range("Q138").Value = 0.25
$ws.Range("R138").Value = 2.025
$ws.Range("S138").Value = 1.825
$ws.Range("U138").Value = 1.975
$ws.Range("V138").Value = 1.875
$ws.Range("W138").Value = -1
$ws.Range("X138").Value = -1
$ws.Range("Y138").Value = 1.05
$ws.Range("Z138").Value = -1
$ws.Range("AA138").Value = 0.825
$ws.Range("AB138").Value = 0.9750000000000001
$ws.Range("AC138").Value = -1
# Row 139
$ws.Range("H139").Value = 1
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = "H"
$ws.Range("W139").Value = 2.1
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 0.7749999999999999
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = -1
$ws.Range("AC139").Value = 0.95
# Row 140
$ws.Range("U140").Value = 2.025
$ws.Range("V140").Value = 1.825
# Row 141
$ws.Range("R141").Value = 1.875
$ws.Range("S141").Value = 1.975
$ws.Range("U141").Value = 2.05
$ws.Range("V141").Value = 1.8
# Row 142
$ws.Range("N142").Value = 8
$ws.Range("O142").Value = 3.8
$ws.Range("P142").Value = 1.45
$ws.Range("R142").Value = 2.025
$ws.Range("S142").Value = 1.825
$ws.Range("T142").Value = 2.5
$ws.Range("U142").Value = 2.05
$ws.Range("V142").Value = 1.8
